# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2210
#   *_new  -> *_FV2304
# then freeze the header row and wrap the data in an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data currently on the sheet.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
$lastCol = $usedRange.Columns.Count()

# Rename the header row (row 1) cells: trailing "_old" -> "_FV2210",
# trailing "_new" -> "_FV2304". Any other header (e.g. "diff") is left as-is.
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val -replace '_old$', '_FV2210'
        $newVal = $newVal -replace '_new$', '_FV2304'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Freeze panes so the header row (row 1) stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the whole sheet range into a native Excel Table with an AutoFilter,
# using the (now renamed) header row as the table's column headers.
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
